$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, which shifts the previous rows 36-44
# down to 37-45 (matching the rest of the diff).
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new record.
$ws.Cells.Item(36, 1).Value = 10
$ws.Cells.Item(36, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(36, 3).Value = "La Araucanía"
$ws.Cells.Item(36, 4).NumberFormat = $ws.Cells.Item(37, 4).NumberFormat
$ws.Cells.Item(36, 4).Value = 44504
$ws.Cells.Item(36, 5).Value = 9
$ws.Cells.Item(36, 6).Value = 100112026
$ws.Cells.Item(36, 7).Value = "Haba"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 155
$ws.Cells.Item(36, 11).Value = 8000
$ws.Cells.Item(36, 12).Value = 9000
$ws.Cells.Item(36, 13).Value = 8613
$ws.Cells.Item(36, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(36, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(36, 16).Value = 345
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"
